$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns used by the claims table (column I / "PolicyId" is intentionally
# left blank for these rows, matching the pattern already used by rows 20-21).
$cols = @("A","B","C","D","E","F","G","H","J","K","L","M","N","O")

# Row 22 duplicates row 2 (NY Windstorm claim), but without the PolicyId (I) value.
foreach ($col in $cols) {
    $ws.Range("$col`2").Copy($ws.Range("$col`22"))
}
$ws.Rows(22).RowHeight = 90

# Rows 23-25 each duplicate row 20 (auto collision claim), again without PolicyId.
foreach ($destRow in 23..25) {
    foreach ($col in $cols) {
        $ws.Range("$col`20").Copy($ws.Range("$col$destRow"))
    }
    $ws.Rows($destRow).RowHeight = 60
}
